# Rename the author of the "FOFA" (SWOT) analysis section heading from
# "Juliana" to "Rafael" in the body text of the document.
#
# Only the exact run text "FOFA Juliana" should be touched - there is an
# unrelated, standalone "Juliana" elsewhere in the document (a tracked
# author/reviewer name) that must be left alone, so we match the full
# "FOFA Juliana" phrase rather than just "Juliana".

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "FOFA Juliana",   # FindText
    $true,            # MatchCase
    $true,            # MatchWholeWord
    $false,           # MatchWildcards
    $false,           # MatchSoundsLike
    $false,           # MatchAllWordForms
    $true,            # Forward
    1,                # Wrap (wdFindContinue)
    $false,           # Format
    "FOFA Rafael",    # ReplaceWith
    2                 # Replace (wdReplaceAll)
)
